$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the multi-row card text blocks into single rows that look like
# Python tuples: (name, [field1, field2, ...])

$hedgeTroll = "('Hedge Troll', ['{2}{G}', 'Creature " + [char]0x2014 + " Troll Cleric', 'Hedge Troll gets +1/+1 as long as you control a Plains.', '{W}: Regenerate Hedge Troll.', '2/2'])"
$oros = "('Oros, the Avenger', ['{3}{R}{W}{B}', 'Legendary Creature " + [char]0x2014 + " Dragon', 'Flying', 'Whenever Oros, the Avenger deals combat damage to a player, you may pay {2}{W}. If you do, Oros deals 3 damage to each nonwhite creature.', '6/6'])"

$ws.Range("A2").Value = $hedgeTroll
$ws.Range("A3").Value = $oros

# Delete the now-unused rows 4 through 13 (original card detail rows)
$ws.Range("A4:A13").EntireRow.Delete()
